$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A58").Value = "Bruno 🐻"
$ws.Range("B58").Value = "Nicolas Giordani  | FC SAVIGNANO"
$ws.Range("C58").Value = "Marco Sala | IMONTAGNA"
$ws.Range("D58").Value = "Danny Giordani | I Magnifici"
$ws.Range("E58").Value = "Riccardo Zaffoni | U.SGUARNA"
$ws.Range("F58").Value = "Alessio Debiasi | Mai una gioia"
